# "Add files via upload" — re-upload of "FINDING MINIMUM OF A FUNCTION.xlsx"
# with the learning-rate (alpha) on the BatchGradientDescent sheet lowered
# from 0.3 to 0.1 (all of the dependent gradient-descent table values
# recalculate from this single input change) plus the last-used-cell
# selections that Excel persists per sheet.

$wb = $excel.ActiveWorkbook

$wsFunction = $wb.Worksheets.Item("FUNCTION")
$wsBatch    = $wb.Worksheets.Item("BatchGradientDescent")

# Lower the learning rate alpha in B2; every formula on this sheet
# (F9:F52, G9:G52, I9:I52, K9:K52, M9:M52, O9:O33) chains off this cell,
# so the recalculation after the script runs reproduces every changed
# cached <v> in the diff.
$wsBatch.Range("B2").Value = 0.1

# Restore the two sheets' remembered selections. Select the
# BatchGradientDescent one first so that selecting on FUNCTION afterwards
# leaves FUNCTION as the active/tabSelected sheet, matching the original.
$wsBatch.Range("S15").Select()
$wsFunction.Range("P36").Select()
